$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Mutual Fund"), shifting existing
# columns C..I to D..J. This matches the dimension growing from A1:I38 to A1:J38.
$ws.Columns.Item(3).Insert()

# Header for the newly inserted column (match the formatting of the other header cells).
$headerSrc = $ws.Cells.Item(1, 4)
$headerDst = $ws.Cells.Item(1, 3)
$headerDst.Value = "Industry"
$headerDst.Font.Bold = $headerSrc.Font.Bold
$headerDst.HorizontalAlignment = $headerSrc.HorizontalAlignment
$headerDst.VerticalAlignment = $headerSrc.VerticalAlignment
$headerDst.Borders.LineStyle = $headerSrc.Borders.LineStyle

# Industry values for each holding row (rows 2-38).
$industries = @(
    "Transport Services",
    "Banks",
    "Healthcare Services",
    "Retailing",
    "Realty",
    "Banks",
    "Pharmaceuticals & Biotechnology",
    "Finance",
    "Retailing",
    "Retailing",
    "Finance",
    "Industrial Manufacturing",
    "Consumer Durables",
    "Finance",
    "Commercial Services & Supplies",
    "Capital Markets",
    "Cement & Cement Products",
    "IT - Software",
    "Healthcare Services",
    "Healthcare Services",
    "Banks",
    "Healthcare Services",
    "Capital Markets",
    "IT - Software",
    "Pharmaceuticals & Biotechnology",
    "Realty",
    "Realty",
    "Electrical Equipment",
    "Consumer Durables",
    "Industrial Products",
    "Consumer Durables",
    "Industrial Products",
    "Industrial Products",
    "Industrial Products",
    "IT - Software",
    "Capital Markets",
    "Beverages"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
